$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Phone Number column holds digit strings, not numbers - force text format
$ws.Range("B2:B3").NumberFormat = "@"

# Row 2: first order
$ws.Range("A2").Value = "Chirayu Sahu"
$ws.Range("B2").Value = "123124234"
$ws.Range("C2").Value = "paracetamol"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "Vit Vellore"

# Row 3: second order
$ws.Range("A3").Value = "Chirayu Sahu"
$ws.Range("B3").Value = "123124234"
$ws.Range("C3").Value = "aspirin"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "Vit Vellore"
